# Regenerate order sheet with updated distance/size codes.
#
# The experiment's distance and size condition labels changed:
#   D64 -> D69
#   D51 -> D55
#   D80 -> D86
#   S30 -> S31   (S20 / S25 are unchanged)
#
# These tokens appear embedded inside many compound strings throughout the
# sheet (Condition names like "Face16_D64_S20", image filenames like
# "Face16_D64_S20_l.png"/"_r.png", and the standalone Distance/Size lookup
# columns such as "D64"/"S30"). Rather than touching individual cells, run
# the substitution as a workbook-wide Find/Replace over the sheet's used
# range so every occurrence (regardless of which column/row it lives in)
# gets updated consistently.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.UsedRange

# Order matters only in that these four tokens are mutually exclusive
# (no token is a substring of another), so a straightforward sequential
# replace is safe.
$range.Replace("D64", "D69") | Out-Null
$range.Replace("D51", "D55") | Out-Null
$range.Replace("D80", "D86") | Out-Null
$range.Replace("S30", "S31") | Out-Null
